$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.189.86"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "2.271.66"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'498.05"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").Value = "'128.99"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("E12").Value = "  +1.48%  "

$ws.Range("D13").Value = "2.674.58"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "'22.71"
$ws.Range("E14").Value = "  +4.91%  "

$ws.Range("D15").Value = "54.132.86"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").Value = "'0.0000130"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "2.281.45"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "'10.21"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("E19").Value = "  +2.06%  "

$ws.Range("D20").Value = "'302.73"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").Value = "'60.98"
$ws.Range("E23").Value = "  -3.61%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").Value = "'0.150"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").Value = "'7.28"
$ws.Range("E26").Value = "  +2.53%  "

$ws.Range("D27").Value = "'170.37"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0684"
$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'5.91"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D33").Value = "'17.73"
$ws.Range("E33").Value = "  +0.70%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "'0.954"
$ws.Range("E35").Value = "  +10.55%  "

$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").Value = "'0.371"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("D41").Value = "'4.80"
$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").Value = "'124.79"
$ws.Range("E42").Value = "  -3.90%  "

$ws.Range("D43").Value = "'0.0491"
$ws.Range("E43").Value = "  +1.95%  "

$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Value = "'0.544"
$ws.Range("E45").Value = "  -0.53%  "

$ws.Range("D46").Value = "'238.70"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("D48").Value = "'0.0205"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").Value = "'16.17"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("E51").Value = "  -0.35%  "
